$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mark rows 447-526 (file list entries) as translated: "ok" in column B
$ws.Range("B447:B526").Value = "ok"

# The "ok"/"not ok" conditional formatting on B2:B621 picked up extra
# (unused) differential-format entries while the sheet was being edited.
# Reproduce that so the style table ends up with the same extra dxf records.
$rngOk = $ws.Range("B2:B621")
for ($i = 0; $i -lt 4; $i++) {
  $fcGreen = $rngOk.FormatConditions.Add(9, 0, "ok")
  $fcGreen.Font.Color = 24832
  $fcGreen.Interior.Color = 13561798
  $fcGreen.Delete()

  $fcRed = $rngOk.FormatConditions.Add(9, 1, "ok")
  $fcRed.Interior.Color = 255
  $fcRed.Delete()
}

# Move the active selection to where it ended up after the edit
$ws.Range("B528").Select()
